# Realestate Update resale numbers 2025-01-11 12:23
# Appends a new data row (row 19) to the CityResaleNum sheet, matching the
# existing data rows (2..18): Date/Time/Weekday/Week as text, city columns
# as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

# Date and Week ("01") look like a date / a number to the engine's literal
# parser, so use a leading apostrophe to force them to be stored as plain
# text (same as the existing rows), then reset the cell style to "Normal"
# so no quote-prefix style/number-format gets stamped onto the cell.
$ws.Cells.Item($row, 1).Value = "'2025-01-11"
$ws.Cells.Item($row, 2).Value = "12:23:56"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 4).Value = "'01"

$ws.Cells.Item($row, 5).Value = 127322
$ws.Cells.Item($row, 6).Value = 143786
$ws.Cells.Item($row, 7).Value = 169434
$ws.Cells.Item($row, 8).Value = 159737
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142888
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192935
$ws.Cells.Item($row, 14).Value = 115441
$ws.Cells.Item($row, 15).Value = 45813
$ws.Cells.Item($row, 16).Value = 28486
$ws.Cells.Item($row, 17).Value = 65112
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48743
$ws.Cells.Item($row, 20).Value = -1

$ws.Range("A19").Style = "Normal"
$ws.Range("D19").Style = "Normal"
